$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.023065755193443
$ws.Range("C2").Value = 0.1927725476481896
$ws.Range("D2").Value = 0.207869663336254
$ws.Range("E2").Value = 0.1706476245223456
$ws.Range("F2").Value = 1.255145777563001
$ws.Range("J2").Value = 0.1839599612594895
$ws.Range("N2").Value = 1.056347490988209
$ws.Range("O2").Value = 2.887423348435817

$ws.Range("B3").Value = 0.9226797810773064
$ws.Range("C3").Value = 0.1685860552133249
$ws.Range("D3").Value = 0.2033822855064784
$ws.Range("E3").Value = 0.1668010292405455
$ws.Range("F3").Value = 1.248864296630501
$ws.Range("J3").Value = 0.1796355010558415
$ws.Range("N3").Value = 1.064277073866833
$ws.Range("O3").Value = 2.881533560150871

$ws.Range("B4").Value = 0.8611533150804576
$ws.Range("C4").Value = 0.1536945318088669
$ws.Range("D4").Value = 0.2007029652278334
$ws.Range("E4").Value = 0.1645237336687053
$ws.Range("F4").Value = 1.245831171726437
$ws.Range("J4").Value = 0.1770892791536269
$ws.Range("N4").Value = 1.069608717379921
$ws.Range("O4").Value = 2.879976776102978

$ws.Range("B5").Value = 0.8361099128395324
$ws.Range("C5").Value = 0.1476161486325509
$ws.Range("D5").Value = 0.199630304737255
$ws.Range("E5").Value = 0.1636169988557015
$ws.Range("F5").Value = 1.244802045488129
$ws.Range("J5").Value = 0.1760790637621525
$ws.Range("N5").Value = 1.071897886162361
$ws.Range("O5").Value = 2.879859350500169

$ws.Range("B6").Value = 0.8319532710856663
$ws.Range("C6").Value = 0.1466062453235679
$ws.Range("D6").Value = 0.1994533515173913
$ws.Range("E6").Value = 0.1634677223633823
$ws.Range("F6").Value = 1.244643649030728
$ws.Range("J6").Value = 0.1759129719728563
$ws.Range("N6").Value = 1.07228503891784
$ws.Range("O6").Value = 2.879871049352658

$ws.Range("B7").Value = 0.8608154514010948
$ws.Range("C7").Value = 0.1536125965313602
$ws.Range("D7").Value = 0.2006884211612459
$ws.Range("E7").Value = 0.1645114189223023
$ws.Range("F7").Value = 1.245816455169162
$ws.Range("J7").Value = 0.1770755441596066
$ws.Range("N7").Value = 1.069639118158605
$ws.Range("O7").Value = 2.879973100468391

$ws.Range("B8").Value = 0.9884303583442033
$ws.Range("C8").Value = 0.1844417240303926
$ws.Range("D8").Value = 0.2063067025765832
$ws.Range("E8").Value = 0.1693037883257276
$ws.Range("F8").Value = 1.252808809370407
$ws.Range("J8").Value = 0.1824462465179408
$ws.Range("N8").Value = 1.058985601662734
$ws.Range("O8").Value = 2.884964521375366

$ws.Range("B9").Value = 1.23952128261061
$ws.Range("C9").Value = 0.2445615781982156
$ws.Range("D9").Value = 0.2179236783230749
$ws.Range("E9").Value = 0.1793718603240748
$ws.Range("F9").Value = 1.273070072568558
$ws.Range("J9").Value = 0.1938449379016163
$ws.Range("N9").Value = 1.041762886983051
$ws.Range("O9").Value = 2.911141077063832

$ws.Range("B10").Value = 1.424471702277174
$ws.Range("C10").Value = 0.288515908351286
$ws.Range("D10").Value = 0.2268211270905738
$ws.Range("E10").Value = 0.1871778862981373
$ws.Range("F10").Value = 1.291970975281203
$ws.Range("J10").Value = 0.202751576907886
$ws.Range("N10").Value = 1.031341643504305
$ws.Range("O10").Value = 2.940435813039528

$ws.Range("B11").Value = 1.508706660822043
$ws.Range("C11").Value = 0.308462990223461
$ws.Range("D11").Value = 0.2309469371558919
$ws.Range("E11").Value = 0.1908180179119228
$ws.Range("F11").Value = 1.301446358644796
$ws.Range("J11").Value = 0.2069198117532665
$ws.Range("N11").Value = 1.02708470926904
$ws.Range("O11").Value = 2.955963794899986

$ws.Range("B12").Value = 1.540617658665553
$ws.Range("C12").Value = 0.3160092671065797
$ws.Range("D12").Value = 0.2325204645426311
$ws.Range("E12").Value = 0.1922092487580187
$ws.Range("F12").Value = 1.305160925925193
$ws.Range("J12").Value = 0.2085150194369163
$ws.Range("N12").Value = 1.02554222902755
$ws.Range("O12").Value = 2.962161591097868

$ws.Range("B13").Value = 1.533744495399219
$ws.Range("C13").Value = 0.3143843694167003
$ws.Range("D13").Value = 0.2321810817752947
$ws.Range("E13").Value = 0.1919090537717381
$ws.Range("F13").Value = 1.304355299006957
$ws.Range("J13").Value = 0.2081707158658901
$ws.Range("N13").Value = 1.025871337917366
$ws.Range("O13").Value = 2.960812638224553

$ws.Range("B14").Value = 1.511331744020083
$ws.Range("C14").Value = 0.309083974285187
$ws.Range("D14").Value = 0.2310761687940754
$ws.Range("E14").Value = 0.1909322190827254
$ws.Range("F14").Value = 1.301749422487234
$ws.Range("J14").Value = 0.2070507138490711
$ws.Range("N14").Value = 1.026956415208005
$ws.Range("O14").Value = 2.956467318163504

$ws.Range("B15").Value = 1.49760494594733
$ws.Range("C15").Value = 0.305836374589461
$ws.Range("D15").Value = 0.2304008300598639
$ws.Range("E15").Value = 0.1903355447879989
$ws.Range("F15").Value = 1.300169723547157
$ws.Range("J15").Value = 0.2063668673021084
$ws.Range("N15").Value = 1.027630110142908
$ws.Range("O15").Value = 2.953847088695142

$ws.Range("B16").Value = 1.418968648510599
$ws.Range("C16").Value = 0.2872113242813157
$ws.Range("D16").Value = 0.2265530639960787
$ws.Range("E16").Value = 0.1869417865490348
$ws.Range("F16").Value = 1.291369414790438
$ws.Range("J16").Value = 0.2024815207006014
$ws.Range("N16").Value = 1.031629575323244
$ws.Range("O16").Value = 2.939465419121035

$ws.Range("B17").Value = 1.370752502622054
$ws.Range("C17").Value = 0.2757729229537631
$ws.Range("D17").Value = 0.2242125792822804
$ws.Range("E17").Value = 0.1848826366529366
$ws.Range("F17").Value = 1.286195609664418
$ws.Range("J17").Value = 0.2001278590129942
$ws.Range("N17").Value = 1.034206990099015
$ws.Range("O17").Value = 2.931207404103304

$ws.Range("B18").Value = 1.34302935130313
$ws.Range("C18").Value = 0.2691893635182225
$ws.Range("D18").Value = 0.2228737711618152
$ws.Range("E18").Value = 0.1837066590382292
$ws.Range("F18").Value = 1.283302328722215
$ws.Range("J18").Value = 0.1987850646698917
$ws.Range("N18").Value = 1.035734984557472
$ws.Range("O18").Value = 2.926664777611421

$ws.Range("B19").Value = 1.333644441410115
$ws.Range("C19").Value = 0.2669595224290049
$ws.Range("D19").Value = 0.2224217434245617
$ws.Range("E19").Value = 0.1833099351380696
$ws.Range("F19").Value = 1.282336883414374
$ws.Range("J19").Value = 0.1983323008140871
$ws.Range("N19").Value = 1.036260158783605
$ws.Range("O19").Value = 2.925162269982309

$ws.Range("B20").Value = 1.375884220153296
$ws.Range("C20").Value = 0.2769910279284602
$ws.Range("D20").Value = 0.2244609650525149
$ws.Range("E20").Value = 0.1851009684252745
$ws.Range("F20").Value = 1.286737823787846
$ws.Range("J20").Value = 0.200377274814187
$ws.Range("N20").Value = 1.033927907589906
$ws.Range("O20").Value = 2.932065035159809

$ws.Range("B21").Value = 1.517914568988544
$ws.Range("C21").Value = 0.3106410283911885
$ws.Range("D21").Value = 0.2314004060063866
$ws.Range("E21").Value = 0.1912187921368016
$ws.Range("F21").Value = 1.302511397545928
$ws.Range("J21").Value = 0.2073792297688328
$ws.Range("N21").Value = 1.026635814959555
$ws.Range("O21").Value = 2.957735013381551

$ws.Range("B22").Value = 1.610815041000251
$ws.Range("C22").Value = 0.332590826921944
$ws.Range("D22").Value = 0.2360008134489675
$ws.Range("E22").Value = 0.1952916936769071
$ws.Range("F22").Value = 1.313557496809665
$ws.Range("J22").Value = 0.2120532754537692
$ws.Range("N22").Value = 1.022275264201454
$ws.Range("O22").Value = 2.976364054784341

$ws.Range("B23").Value = 1.561225851994607
$ws.Range("C23").Value = 0.3208798042343517
$ws.Range("D23").Value = 0.2335395638135935
$ws.Range("E23").Value = 0.1931110963544498
$ws.Range("F23").Value = 1.307594438182406
$ws.Range("J23").Value = 0.2095496860554533
$ws.Range("N23").Value = 1.024565500921391
$ws.Range("O23").Value = 2.966251545344676

$ws.Range("B24").Value = 1.373564180044355
$ws.Range("C24").Value = 0.2764403459146081
$ws.Range("D24").Value = 0.2243486487598148
$ws.Range("E24").Value = 0.185002236162461
$ws.Range("F24").Value = 1.286492435841666
$ws.Range("J24").Value = 0.2002644816732442
$ws.Range("N24").Value = 1.03405393686144
$ws.Range("O24").Value = 2.93167666161537

$ws.Range("B25").Value = 1.171508640566856
$ws.Range("C25").Value = 0.2283346660692587
$ws.Range("D25").Value = 0.2147171087037236
$ws.Range("E25").Value = 0.1765763831437752
$ws.Range("F25").Value = 1.266885412408286
$ws.Range("J25").Value = 0.1906681233375451
$ws.Range("N25").Value = 1.046029786764343
$ws.Range("O25").Value = 2.902297696862945

